$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Recorded By" (column G) values for each session row to the
# academic year(s) the session was recorded in. Row 29 has no "Recorded By"
# value (session not yet recorded), so it is left untouched.
$rows = 2..30
foreach ($r in $rows) {
    if ($r -eq 29) { continue }
    if ($r -eq 27) {
        $ws.Range("G$r").Value = "2022/2023, 2023/2024, 2025/2026"
    } else {
        $ws.Range("G$r").Value = "2025/2026"
    }
}

# Narrow column G now that it holds short academic-year strings instead of
# long lists of recorder names. Excel's ColumnWidth property is offset from
# the stored OOXML column width by the default column padding (~0.8333
# characters), so back that out to land on a stored width of exactly 33.
$ws.Columns.Item(7).ColumnWidth = 33 - 0.8333333333333334
